$wb = $excel.ActiveWorkbook

# --- BoM sheet ---
$bom = $wb.Worksheets.Item("BoM")

# Row 16: Resistor group R2 R4 -> R2 only (qty 2 -> 1)
$bom.Range("D16").Value = "R2"
$bom.Range("G16").Value = 1

# Row 18: U1 U2 -> U1 only (qty 2 -> 1)
$bom.Range("D18").Value = "U1"
$bom.Range("G18").Value = 1

# Summary fields
$bom.Range("F3").Value = "216 (205 SMD/ 0 THT)"
$bom.Range("F4").Value = "213 (204 SMD/ 0 THT)"
$bom.Range("F6").Value = 213

# --- DNF sheet ---
$dnf = $wb.Worksheets.Item("DNF")
$dnf.Range("F3").Value = "216 (205 SMD/ 0 THT)"
$dnf.Range("F4").Value = "213 (204 SMD/ 0 THT)"
$dnf.Range("F6").Value = 213
